$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.925.06"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "3.439.98"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").Value = "'578.02"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "'146.56"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").Value = "3.452.39"
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.477"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").Value = "'7.78"
$ws.Range("E10").Value = "  +2.37%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "'0.404"
$ws.Range("E12").Value = "  +3.95%  "
$ws.Range("D13").Value = "4.050.30"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").Value = "'0.128"
$ws.Range("E14").Value = "  +2.33%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'28.97"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").Value = "3.454.20"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").Value = "'0.0000171"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "63.043.74"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").Value = "'6.39"
$ws.Range("E19").Value = "  +3.06%  "
$ws.Range("D20").Value = "'14.38"
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("D21").Value = "'9.20"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "'383.94"
$ws.Range("E22").Value = "  -2.04%  "
$ws.Range("D23").Value = "'0.560"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("D24").Value = "'74.41"
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "3.601.05"
$ws.Range("E26").Value = "  +1.24%  "
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").Value = "'0.180"
$ws.Range("E28").Value = "  -4.12%  "
$ws.Range("D29").Value = "'7.60"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "'8.11"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").Value = "'2.11"
$ws.Range("E32").Value = "  -1.48%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'23.23"
$ws.Range("E34").Value = "  -1.91%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "'1.33"
$ws.Range("E35").Value = "  -4.81%  "
$ws.Range("D36").Value = "'5.29"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("D37").Value = "'7.09"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'1.59"
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "'31.80"
$ws.Range("E39").Value = "  +9.20%  "
$ws.Range("D40").Value = "'168.32"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "3.493.23"
$ws.Range("E41").Value = "  +1.16%  "
$ws.Range("D42").Value = "'0.0767"
$ws.Range("E42").Value = "  +1.23%  "
$ws.Range("D43").Value = "'0.789"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "'42.30"
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("D45").Value = "'1.72"
$ws.Range("E45").Value = "  +1.75%  "
$ws.Range("D46").Value = "'1.20"
$ws.Range("E46").Value = "  +2.37%  "
$ws.Range("D47").Value = "'4.34"
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("D48").Value = "2.585.61"
$ws.Range("E48").Value = "  +3.11%  "
$ws.Range("D49").Value = "'2.28"
$ws.Range("E49").Value = "  +7.86%  "
$ws.Range("D50").Value = "'6.77"
$ws.Range("E50").Value = "  +1.51%  "
$ws.Range("D51").Value = "'22.78"
$ws.Range("E51").Value = "  -0.95%  "
